$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '44.467.40'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.250.58'
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '308.49'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '94.86'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -2.42%  '
$c.Style = "Normal"

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.572'
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -0.26%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.09'
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +0.10%  '
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +1.00%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '2.362.83'
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +3.07%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.844'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +1.86%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '13.73'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '44.199.56'
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.57%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.0₃0968'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.48'
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -0.38%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '6.42'
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +2.55%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '65.92'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +1.36%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '238.01'
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -0.37%  '
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +3.02%  '
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +3.94%  '
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -0.24%  '
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '38.58'
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +6.18%  '
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +4.40%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '9.91'
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.Style = "Normal"

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '6.00'
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -1.22%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '20.13'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +0.54%  '
$c.Style = "Normal"

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '154.46'
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -0.33%  '
$c.Style = "Normal"

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.0803'
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -0.79%  '
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '3.13'
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -8.66%  '
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +3.29%  '
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +1.14%  '
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  +0.41%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '3.49'
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +5.14%  '
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '14.84'
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -1.55%  '
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +0.39%  '
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.743.76'
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +2.31%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '80.94'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -6.02%  '
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '71.40'
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +4.79%  '
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '100.26'
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '4.97'
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -3.05%  '
$c.Style = "Normal"

$c = $ws.Range('B49')
$c.NumberFormat = "@"
$c.Value = 'MultiversX'
$c.Style = "Normal"

$c = $ws.Range('C49')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '56.39'
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +2.78%  '
$c.Style = "Normal"

$c = $ws.Range('B50')
$c.NumberFormat = "@"
$c.Value = 'Stacks'
$c.Style = "Normal"

$c = $ws.Range('C50')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.61'
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +6.77%  '
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '8.14'
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -0.65%  '
$c.Style = "Normal"

